$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Digits (keycap emoji, e.g. "2" + VARIATION SELECTOR-16 + COMBINING ENCLOSING KEYCAP)
# for column A ("Раздел") on rows 4..24, grouping the questions into sections.
$fe0f = [char]0xFE0F
$e20e3 = [char]0x20E3

$rowDigits = @{
    4  = 2
    5  = 2
    6  = 3
    7  = 3
    8  = 3
    9  = 4
    10 = 4
    11 = 5
    12 = 5
    13 = 6
    14 = 6
    15 = 6
    16 = 7
    17 = 7
    18 = 7
    19 = 8
    20 = 8
    21 = 8
    22 = 9
    23 = 9
    24 = 9
}

foreach ($row in $rowDigits.Keys) {
    $digit = $rowDigits[$row]
    $value = "$digit$fe0f$e20e3"
    $ws.Cells.Item($row, 1).Value = $value
}

[void]$ws.Range("B2").Select()
